# Refresh the cryptocurrency price/volume snapshot (GitHub Actions bot update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "42.867.01" },
    @{ Cell = "E2"; Value = "  -5.79%  " },
    @{ Cell = "D3"; Value = "2.216.26" },
    @{ Cell = "E3"; Value = "  -6.87%  " },
    @{ Cell = "E4"; Value = "  -0.09%  " },
    @{ Cell = "D5"; Value = "'313.00" },
    @{ Cell = "E5"; Value = "  -0.42%  " },
    @{ Cell = "D6"; Value = "'99.41" },
    @{ Cell = "E6"; Value = "  -8.61%  " },
    @{ Cell = "D7"; Value = "'0.581" },
    @{ Cell = "E7"; Value = "  -7.81%  " },
    @{ Cell = "D8"; Value = "'0.999" },
    @{ Cell = "E8"; Value = "  -0.15%  " },
    @{ Cell = "D9"; Value = "'0.558" },
    @{ Cell = "E9"; Value = "  -9.20%  " },
    @{ Cell = "D10"; Value = "'36.63" },
    @{ Cell = "E10"; Value = "  -10.66%  " },
    @{ Cell = "D11"; Value = "'54.46" },
    @{ Cell = "E11"; Value = "  -2.99%  " },
    @{ Cell = "D12"; Value = "'0.0821" },
    @{ Cell = "E12"; Value = "  -10.82%  " },
    @{ Cell = "D13"; Value = "'7.59" },
    @{ Cell = "E13"; Value = "  -11.34%  " },
    @{ Cell = "E14"; Value = "  -1.21%  " },
    @{ Cell = "D15"; Value = "2.549.29" },
    @{ Cell = "E15"; Value = "  -6.89%  " },
    @{ Cell = "D16"; Value = "'0.850" },
    @{ Cell = "E16"; Value = "  -13.86%  " },
    @{ Cell = "D17"; Value = "'14.10" },
    @{ Cell = "E17"; Value = "  -8.41%  " },
    @{ Cell = "D18"; Value = "2.211.85" },
    @{ Cell = "E18"; Value = "  -6.36%  " },
    @{ Cell = "D19"; Value = "42.748.62" },
    @{ Cell = "E19"; Value = "  -6.01%  " },
    @{ Cell = "D20"; Value = "'14.35" },
    @{ Cell = "E20"; Value = "  +3.36%  " },
    @{ Cell = "D21"; Value = "0.0₃0958" },
    @{ Cell = "E21"; Value = "  -10.26%  " },
    @{ Cell = "D22"; Value = "'6.46" },
    @{ Cell = "E22"; Value = "  -11.25%  " },
    @{ Cell = "D23"; Value = "'65.16" },
    @{ Cell = "E23"; Value = "  -11.43%  " },
    @{ Cell = "D24"; Value = "'3.10" },
    @{ Cell = "E24"; Value = "  -12.52%  " },
    @{ Cell = "D25"; Value = "'236.61" },
    @{ Cell = "E25"; Value = "  -9.27%  " },
    @{ Cell = "E26"; Value = "  -11.04%  " },
    @{ Cell = "E27"; Value = "  -0.26%  " },
    @{ Cell = "D28"; Value = "'4.06" },
    @{ Cell = "E28"; Value = "  +1.60%  " },
    @{ Cell = "E29"; Value = "  -3.60%  " },
    @{ Cell = "D30"; Value = "'9.90" },
    @{ Cell = "E30"; Value = "  -11.31%  " },
    @{ Cell = "D31"; Value = "'6.33" },
    @{ Cell = "E31"; Value = "  -13.11%  " },
    @{ Cell = "D32"; Value = "'20.34" },
    @{ Cell = "E32"; Value = "  -9.32%  " },
    @{ Cell = "D33"; Value = "'0.0868" },
    @{ Cell = "E33"; Value = "  -11.86%  " },
    @{ Cell = "D34"; Value = "'34.07" },
    @{ Cell = "E34"; Value = "  -8.74%  " },
    @{ Cell = "D35"; Value = "'155.30" },
    @{ Cell = "E35"; Value = "  -6.79%  " },
    @{ Cell = "D36"; Value = "'2.76" },
    @{ Cell = "E36"; Value = "  -7.16%  " },
    @{ Cell = "B37"; Value = "Stellar" },
    @{ Cell = "C37"; Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm" },
    @{ Cell = "D37"; Value = "'0.121" },
    @{ Cell = "E37"; Value = "  -7.53%  " },
    @{ Cell = "B38"; Value = "LidoDAOToken" },
    @{ Cell = "C38"; Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo" },
    @{ Cell = "D38"; Value = "'2.97" },
    @{ Cell = "E38"; Value = "  -0.63%  " },
    @{ Cell = "D39"; Value = "'1.88" },
    @{ Cell = "E39"; Value = "  -3.13%  " },
    @{ Cell = "E40"; Value = "  -7.19%  " },
    @{ Cell = "D41"; Value = "'0.105" },
    @{ Cell = "E41"; Value = "  -12.03%  " },
    @{ Cell = "D42"; Value = "'3.67" },
    @{ Cell = "E42"; Value = "  -8.94%  " },
    @{ Cell = "D43"; Value = "'0.0321" },
    @{ Cell = "E43"; Value = "  -10.27%  " },
    @{ Cell = "E44"; Value = "  -0.05%  " },
    @{ Cell = "D45"; Value = "1.775.13" },
    @{ Cell = "E45"; Value = "  -2.61%  " },
    @{ Cell = "D46"; Value = "'12.21" },
    @{ Cell = "E46"; Value = "  -4.88%  " },
    @{ Cell = "B47"; Value = "BitcoinSV" },
    @{ Cell = "C47"; Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv" },
    @{ Cell = "D47"; Value = "'86.03" },
    @{ Cell = "E47"; Value = "  -13.22%  " },
    @{ Cell = "B48"; Value = "Algorand" },
    @{ Cell = "C48"; Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo" },
    @{ Cell = "D48"; Value = "'0.203" },
    @{ Cell = "E48"; Value = "  -11.50%  " },
    @{ Cell = "D49"; Value = "'5.29" },
    @{ Cell = "E49"; Value = "  -9.14%  " },
    @{ Cell = "D50"; Value = "'75.33" },
    @{ Cell = "E50"; Value = "  -10.58%  " },
    @{ Cell = "D51"; Value = "'58.83" },
    @{ Cell = "E51"; Value = "  -15.71%  " }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
